$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.296680569648743
$ws.Range("B1").Value = 2.296324968338013
$ws.Range("C1").Value = 2.917322158813477
$ws.Range("D1").Value = 3.35934591293335
$ws.Range("E1").Value = 1.699271678924561
